$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 currently holds the "law-salary" topic, and row 30 is an empty
# placeholder row (it only carries the wrap-text formatting on column B).
# Insert a fresh row above row 29 so the existing "law-salary" content
# shifts down into row 30 (taking its formatting/height with it), then
# drop the old placeholder row that is now duplicated further down so the
# sheet keeps the same number of rows.
$ws.Rows.Item(29).Insert()
$ws.Rows.Item(31).Delete()

# Give the new row 29 the same formatting as row 30 (which still has the
# original look of this two-column topic/response table).
$ws.Range("B30").Copy()
$ws.Range("B29").PasteSpecial(-4122)
$ws.Range("A30").Copy()
$ws.Range("A29").PasteSpecial(-4122)

# Add the new "law-place" topic and its response text.
$ws.Range("A29").Value = "law-place"
$ws.Range("B29").Value = "1.สำนักงานสรรพากรทุกสาขาทุกเเห่ง `n2.ไปรษณีย์ เเบบลงทะเบียน `n3.ช่องทางออนไลน์ ผ่านเว็บไซต์ของกรมสรรพากร"
$ws.Rows.Item(29).RowHeight = 38.25

# Reflect the new scroll position / selection used while editing.
$excel.ActiveWindow.ScrollRow = 25
$ws.Range("B30").Select()
